$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.277433666666667
$ws.Range("H2").Value = 3.832301
$ws.Range("I2").Value = 0.01786062203930835
$ws.Range("J2").Value = 0.01786062203930835
$ws.Range("M2").Value = 247.0944516666667
$ws.Range("N2").Value = 741.283355
$ws.Range("O2").Value = 0.8050739182622993
$ws.Range("P2").Value = 0.8050739182622993
$ws.Range("Q2").Value = 315.6467714055394
$ws.Range("R2").Value = 2840.820942649855
$ws.Range("S2").Value = 0.01437912096778796
$ws.Range("T2").Value = 0.01437912096778795
$ws.Range("G3").Value = 1.277433666666667
$ws.Range("H3").Value = 3.832301
$ws.Range("I3").Value = 0.01786062203930835
$ws.Range("J3").Value = 0.01786062203930835
$ws.Range("O3").Value = 0.1379009747488701
$ws.Range("P3").Value = 0.13790097474887
$ws.Range("Q3").Value = 54.06708187381111
$ws.Range("R3").Value = 486.6037368642999
$ws.Range("S3").Value = 0.002462997188841774
$ws.Range("T3").Value = 0.002462997188841773
$ws.Range("G4").Value = 1.277433666666667
$ws.Range("H4").Value = 3.832301
$ws.Range("I4").Value = 0.01786062203930835
$ws.Range("J4").Value = 0.01786062203930835
$ws.Range("M4").Value = 11.590146
$ws.Range("N4").Value = 34.770438
$ws.Range("O4").Value = 0.03776258103132013
$ws.Range("P4").Value = 0.03776258103132013
$ws.Range("Q4").Value = 14.805642701982
$ws.Range("R4").Value = 133.250784317838
$ws.Range("S4").Value = 0.0006744631870291639
$ws.Range("T4").Value = 0.0006744631870291638
$ws.Range("G5").Value = 1.277433666666667
$ws.Range("H5").Value = 3.832301
$ws.Range("I5").Value = 0.01786062203930835
$ws.Range("J5").Value = 0.01786062203930835
$ws.Range("M5").Value = 5.912082333333333
$ws.Range("N5").Value = 17.736247
$ws.Range("O5").Value = 0.01926252595751047
$ws.Range("P5").Value = 0.01926252595751047
$ws.Range("Q5").Value = 7.552293012705221
$ws.Range("R5").Value = 67.97063711434699
$ws.Range("S5").Value = 0.0003440406956494608
$ws.Range("T5").Value = 0.0003440406956494608
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("I6").Value = 0.6354272679079697
$ws.Range("J6").Value = 0.6354272679079697
$ws.Range("M6").Value = 247.0944516666667
$ws.Range("N6").Value = 741.283355
$ws.Range("O6").Value = 0.8050739182622993
$ws.Range("P6").Value = 0.8050739182622993
$ws.Range("Q6").Value = 11229.763730332
$ws.Range("R6").Value = 101067.873572988
$ws.Range("S6").Value = 0.511565920345377
$ws.Range("T6").Value = 0.511565920345377
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("I7").Value = 0.6354272679079697
$ws.Range("J7").Value = 0.6354272679079697
$ws.Range("O7").Value = 0.1379009747488701
$ws.Range("P7").Value = 0.13790097474887
$ws.Range("S7").Value = 0.08762603962652041
$ws.Range("T7").Value = 0.08762603962652039
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("I8").Value = 0.6354272679079697
$ws.Range("J8").Value = 0.6354272679079697
$ws.Range("M8").Value = 11.590146
$ws.Range("N8").Value = 34.770438
$ws.Range("O8").Value = 0.03776258103132013
$ws.Range("P8").Value = 0.03776258103132013
$ws.Range("Q8").Value = 526.74028211541
$ws.Range("R8").Value = 4740.66253903869
$ws.Range("S8").Value = 0.02399537369388507
$ws.Range("T8").Value = 0.02399537369388507
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("I9").Value = 0.6354272679079697
$ws.Range("J9").Value = 0.6354272679079697
$ws.Range("M9").Value = 5.912082333333333
$ws.Range("N9").Value = 17.736247
$ws.Range("O9").Value = 0.01926252595751047
$ws.Range("P9").Value = 0.01926252595751047
$ws.Range("Q9").Value = 268.6878936770539
$ws.Range("R9").Value = 2418.191043093485
$ws.Range("S9").Value = 0.01223993424218723
$ws.Range("T9").Value = 0.01223993424218723
$ws.Range("G10").Value = 23.96074166666667
$ws.Range("H10").Value = 71.88222500000001
$ws.Range("I10").Value = 0.3350105464235513
$ws.Range("J10").Value = 0.3350105464235513
$ws.Range("M10").Value = 247.0944516666667
$ws.Range("N10").Value = 741.283355
$ws.Range("O10").Value = 0.8050739182622993
$ws.Range("P10").Value = 0.8050739182622993
$ws.Range("Q10").Value = 5920.566323651653
$ws.Range("R10").Value = 53285.09691286488
$ws.Range("S10").Value = 0.2697082532684024
$ws.Range("T10").Value = 0.2697082532684024
$ws.Range("G11").Value = 23.96074166666667
$ws.Range("H11").Value = 71.88222500000001
$ws.Range("I11").Value = 0.3350105464235513
$ws.Range("J11").Value = 0.3350105464235513
$ws.Range("O11").Value = 0.1379009747488701
$ws.Range("P11").Value = 0.13790097474887
$ws.Range("Q11").Value = 1014.132800201945
$ws.Range("R11").Value = 9127.1952018175
$ws.Range("S11").Value = 0.04619828090295931
$ws.Range("T11").Value = 0.04619828090295929
$ws.Range("G12").Value = 23.96074166666667
$ws.Range("H12").Value = 71.88222500000001
$ws.Range("I12").Value = 0.3350105464235513
$ws.Range("J12").Value = 0.3350105464235513
$ws.Range("M12").Value = 11.590146
$ws.Range("N12").Value = 34.770438
$ws.Range("O12").Value = 0.03776258103132013
$ws.Range("P12").Value = 0.03776258103132013
$ws.Range("Q12").Value = 277.70849418495
$ws.Range("R12").Value = 2499.37644766455
$ws.Range("S12").Value = 0.01265086290566619
$ws.Range("T12").Value = 0.01265086290566619
$ws.Range("G13").Value = 23.96074166666667
$ws.Range("H13").Value = 71.88222500000001
$ws.Range("I13").Value = 0.3350105464235513
$ws.Range("J13").Value = 0.3350105464235513
$ws.Range("M13").Value = 5.912082333333333
$ws.Range("N13").Value = 17.736247
$ws.Range("O13").Value = 0.01926252595751047
$ws.Range("P13").Value = 0.01926252595751047
$ws.Range("Q13").Value = 141.6578775010639
$ws.Range("R13").Value = 1274.920897509575
$ws.Range("S13").Value = 0.006453149346523425
$ws.Range("T13").Value = 0.006453149346523424
$ws.Range("G14").Value = 0.8369233333333334
$ws.Range("H14").Value = 2.51077
$ws.Range("I14").Value = 0.01170156362917063
$ws.Range("J14").Value = 0.01170156362917063
$ws.Range("M14").Value = 247.0944516666667
$ws.Range("N14").Value = 741.283355
$ws.Range("O14").Value = 0.8050739182622993
$ws.Range("P14").Value = 0.8050739182622993
$ws.Range("Q14").Value = 206.7991121370389
$ws.Range("R14").Value = 1861.19200923335
$ws.Range("S14").Value = 0.009420623680732012
$ws.Range("T14").Value = 0.00942062368073201
$ws.Range("G15").Value = 0.8369233333333334
$ws.Range("H15").Value = 2.51077
$ws.Range("I15").Value = 0.01170156362917063
$ws.Range("J15").Value = 0.01170156362917063
$ws.Range("O15").Value = 0.1379009747488701
$ws.Range("P15").Value = 0.13790097474887
$ws.Range("Q15").Value = 35.42258480122222
$ws.Range("R15").Value = 318.803263211
$ws.Range("S15").Value = 0.001613657030548556
$ws.Range("T15").Value = 0.001613657030548555
$ws.Range("G16").Value = 0.8369233333333334
$ws.Range("H16").Value = 2.51077
$ws.Range("I16").Value = 0.01170156362917063
$ws.Range("J16").Value = 0.01170156362917063
$ws.Range("M16").Value = 11.590146
$ws.Range("N16").Value = 34.770438
$ws.Range("O16").Value = 0.03776258103132013
$ws.Range("P16").Value = 0.03776258103132013
$ws.Range("Q16").Value = 9.70006362414
$ws.Range("R16").Value = 87.30057261725999
$ws.Range("S16").Value = 0.0004418812447397044
$ws.Range("T16").Value = 0.0004418812447397044
$ws.Range("G17").Value = 0.8369233333333334
$ws.Range("H17").Value = 2.51077
$ws.Range("I17").Value = 0.01170156362917063
$ws.Range("J17").Value = 0.01170156362917063
$ws.Range("M17").Value = 5.912082333333333
$ws.Range("N17").Value = 17.736247
$ws.Range("O17").Value = 0.01926252595751047
$ws.Range("P17").Value = 0.01926252595751047
$ws.Range("Q17").Value = 4.947959653354444
$ws.Range("R17").Value = 44.53163688018999
$ws.Range("S17").Value = 0.0002254016731503597
$ws.Range("T17").Value = 0.0002254016731503597
